$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Jam - Raspberry (Baking)): Quantity 3 -> 7, Total Cost $504.93 -> $1,178.17
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "7"
$ws.Range("C3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "$1,178.17"
$ws.Range("E3").Style = "Normal"

# Row 11 (Flour - Millers Choice): Quantity 100 -> 60, Total Cost $1,689.00 -> $1,013.40
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "60"
$ws.Range("C11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "$1,013.40"
$ws.Range("E11").Style = "Normal"
